$wb = $excel.ActiveWorkbook

$test = $wb.Worksheets.Add()
$test.Name = "Test"

$testRef = $wb.Worksheets.Item("Test")
$namingRef = $wb.Worksheets.Item("Naming")

$testRef.Move([System.Reflection.Missing]::Value, $namingRef)

# re-fetch after move, in case reference went stale
$testRef2 = $wb.Worksheets.Item("Test")
$testRef2.Activate()

Write-Output "After Move(After=naming) + activate:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
Write-Output ("ActiveSheet: " + $wb.ActiveSheet.Name)
